$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin (multi-dot string stays text naturally)
$ws.Range("D2").Value = "38.188.84"
$ws.Range("E2").Value = "  +1.95%  "

# Row 3 - Ethereum (multi-dot string stays text naturally)
$ws.Range("D3").Value = "2.056.38"
$ws.Range("E3").Value = "  +1.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB (looks numeric - force text)
Set-TextValue "D5" "228.32"
$ws.Range("E5").Value = "  -0.62%  "

# Row 6 - XRP (looks numeric - force text)
Set-TextValue "D6" "0.617"
$ws.Range("E6").Value = "  +0.53%  "

# Row 7 - Solana (looks numeric - force text)
Set-TextValue "D7" "60.90"

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.18%  "

# Row 10 - Dogecoin (looks numeric - force text)
Set-TextValue "D10" "0.0827"
$ws.Range("E10").Value = "  +5.19%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.47%  "

# Row 12 - Chainlink (looks numeric - force text)
Set-TextValue "D12" "14.82"
$ws.Range("E12").Value = "  +2.77%  "

# Row 13 - WrappedliquidstakedEther2.0 (multi-dot string stays text naturally)
$ws.Range("D13").Value = "2.358.91"
$ws.Range("E13").Value = "  +1.15%  "

# Row 14 - Avalanche (looks numeric - force text)
Set-TextValue "D14" "21.11"
$ws.Range("E14").Value = "  +3.21%  "

# Row 15 - Polygon (looks numeric - force text)
Set-TextValue "D15" "0.762"
$ws.Range("E15").Value = "  +2.41%  "

# Row 16 - Polkadot (looks numeric - force text)
Set-TextValue "D16" "5.31"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17 - WrappedEther (multi-dot string stays text naturally)
$ws.Range("D17").Value = "2.069.40"
$ws.Range("E17").Value = "  +2.12%  "

# Row 18 - WrappedBTC (multi-dot string stays text naturally)
$ws.Range("D18").Value = "38.121.40"
$ws.Range("E18").Value = "  +2.07%  "

# Row 19 - Uniswap (looks numeric - force text)
Set-TextValue "D19" "6.17"
$ws.Range("E19").Value = "  -1.16%  "

# Row 20 - Litecoin (looks numeric - force text)
Set-TextValue "D20" "69.82"
$ws.Range("E20").Value = "  +1.02%  "

# Row 21 - ShibaInu (has subscript char, stays text naturally)
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +1.41%  "

# Row 22 - BitcoinCash (looks numeric - force text)
Set-TextValue "D22" "225.12"
$ws.Range("E22").Value = "  +0.64%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.12%  "

# Row 24 - Toncoin (looks numeric - force text)
Set-TextValue "D24" "2.43"
$ws.Range("E24").Value = "  -1.00%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.53%  "

# Row 26 - Monero (looks numeric - force text)
Set-TextValue "D26" "166.71"
$ws.Range("E26").Value = "  +1.25%  "

# Row 27 - Cosmos (looks numeric - force text)
Set-TextValue "D27" "9.22"
$ws.Range("E27").Value = "  +0.52%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -1.06%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +1.00%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -2.01%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  +3.01%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.20%  "

# Row 33 - WEMIXToken (looks numeric - force text)
Set-TextValue "D33" "2.06"
$ws.Range("E33").Value = "  +1.08%  "

# Row 34 - InternetComputer(DFINITY) (looks numeric - force text)
Set-TextValue "D34" "4.56"
$ws.Range("E34").Value = "  +2.06%  "

# Row 35 - Hedera (looks numeric - force text)
Set-TextValue "D35" "0.0607"
$ws.Range("E35").Value = "  -0.23%  "

# Row 36 - THORChain (looks numeric - force text)
Set-TextValue "D36" "6.30"
$ws.Range("E36").Value = "  +10.80%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -2.67%  "

# Row 38 - RenderToken (looks numeric - force text)
Set-TextValue "D38" "3.27"
$ws.Range("E38").Value = "  +1.54%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  +0.13%  "

# Row 40 - Maker (multi-dot string stays text naturally)
$ws.Range("D40").Value = "1.532.86"
$ws.Range("E40").Value = "  +4.29%  "

# Row 41 - was VeChain, now Aave (looks numeric - force text)
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D41" "98.35"
$ws.Range("E41").Value = "  +3.89%  "

# Row 42 - was Aave, now VeChain (looks numeric - force text)
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.0219"
$ws.Range("E42").Value = "  +2.07%  "

# Row 43 - InjectiveProtocol (looks numeric - force text)
Set-TextValue "D43" "16.82"
$ws.Range("E43").Value = "  +2.92%  "

# Row 44 - HuobiToken (looks numeric - force text)
Set-TextValue "D44" "2.84"
$ws.Range("E44").Value = "  +0.81%  "

# Row 45 - Cronos (looks numeric - force text)
Set-TextValue "D45" "0.0930"
$ws.Range("E45").Value = "  +0.51%  "

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = "  +1.07%  "

# Row 47 - FTXToken (looks numeric - force text)
Set-TextValue "D47" "3.99"
$ws.Range("E47").Value = "  -8.42%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +0.35%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  +1.67%  "

# Row 50 - FraxShare (looks numeric - force text)
Set-TextValue "D50" "7.06"
$ws.Range("E50").Value = "  -0.89%  "

# Row 51 - RocketPoolETH (multi-dot string stays text naturally)
$ws.Range("D51").Value = "2.247.89"
$ws.Range("E51").Value = "  +1.31%  "
